$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new "Wins" / "Losses" / "Ties" columns (AD1:AF1) ---
# Match the existing header formatting (bold font, thin border, centered/top aligned)
# by copying the format from an existing header cell (A1) onto the new ones.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows 2-46: season record (Wins=99, Losses=62, Ties=0) for every player row ---
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 99
    $ws.Cells.Item($r, 31).Value = 62
    $ws.Cells.Item($r, 32).Value = 0
}
